$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ E=3; G=2.3082;            H=6.9246;             I=0.6638288620319053; J=0.6638288620319053; K=3; M=34.07074633333333; N=102.212239;        O=0.5171464495142372; P=0.5171464495142373; Q=78.64209668659997;  R=707.7788701793999;  S=0.3432967390848762; T=0.3432967390848763 }
    3 = @{ E=3; G=2.3082;            H=6.9246;             I=0.6638288620319053; J=0.6638288620319053; K=3; M=27.685497;          N=83.05649099999999; O=0.420227262899125;  P=0.4202272628991251; Q=63.90366417539999;  R=575.1329775785999;  S=0.2789589857251085; T=0.2789589857251085 }
    4 = @{ E=3; G=2.3082;            H=6.9246;             I=0.6638288620319053; J=0.6638288620319053; K=3; M=4.125957666666666;  N=12.377873;          O=0.06262628758663766;P=0.06262628758663766;Q=9.523535486199998;   R=85.71181937579999;  S=0.04157313722192051;T=0.04157313722192051 }
    5 = @{ E=3; G=1.168901;          H=3.506703;           I=0.3361711379680947; J=0.3361711379680947; K=3; M=34.07074633333333; N=102.212239;        O=0.5171464495142372; P=0.5171464495142373; Q=39.82532945977965;  R=358.4279651380169;  S=0.173849710429361;  T=0.173849710429361 }
    6 = @{ E=3; G=1.168901;          H=3.506703;           I=0.3361711379680947; J=0.3361711379680947; K=3; M=27.685497;          N=83.05649099999999; O=0.420227262899125;  P=0.4202272628991251; Q=32.361605128797;    R=291.254446159173;   S=0.1412682771740166; T=0.1412682771740166 }
    7 = @{ E=3; G=1.168901;          H=3.506703;           I=0.3361711379680947; J=0.3361711379680947; K=3; M=4.125957666666666;  N=12.377873;          O=0.06262628758663766;P=0.06262628758663766;Q=4.822836042524333;   R=43.40552438271899;  S=0.02105315036471715;T=0.02105315036471714 }
}

foreach ($r in $data.Keys) {
    $row = $data[$r]
    foreach ($col in $row.Keys) {
        $ws.Range("$col$r").Value = $row[$col]
    }
}
